$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1108.9
$ws.Range("I4").Value = 1108.9
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1108.9
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -994.9000000000001
$ws.Range("N4").ClearContents()

$ws.Range("H9").Value = 118.31579
$ws.Range("I9").Value = 117.588234
$ws.Range("K9").Value = 117.588234
$ws.Range("M9").Value = 51.411766

$ws.Range("H12").Value = 169
$ws.Range("I12").Value = 221
$ws.Range("J12").Value = 65
$ws.Range("K12").Value = 221
$ws.Range("L12").Value = 65
$ws.Range("M12").Value = -51
$ws.Range("N12").Value = -405

$ws.Range("H17").Value = 1197.3158
$ws.Range("J17").Value = 1230.2941
$ws.Range("L17").Value = 3690.8823
$ws.Range("N17").Value = -4026.8823

$ws.Range("H18").Value = 3885.4614
$ws.Range("I18").Value = 4892.1
$ws.Range("J18").Value = 530
$ws.Range("K18").Value = 4892.1
$ws.Range("L18").Value = 530
$ws.Range("M18").Value = -4608.1
$ws.Range("N18").Value = -1098

$ws.Range("H28").Value = 2512.5386
$ws.Range("J28").Value = 5000
$ws.Range("L28").Value = 5000
$ws.Range("N28").Value = -5970

$ws.Range("H32").Value = 1103.75
$ws.Range("J32").Value = 1171.7858
$ws.Range("L32").Value = 1171.7858
$ws.Range("N32").Value = -1823.7858

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H58").Value = 3251.8333
$ws.Range("J58").Value = 4700.727
$ws.Range("L58").Value = 14102.181
$ws.Range("N58").Value = -14402.181

$ws.Range("H62").Value = 111113240
$ws.Range("I62").Value = 111113240
$ws.Range("K62").Value = 111113240
$ws.Range("M62").Value = -111112616

$ws.Range("H64").Value = 92248.625
$ws.Range("I64").Value = 131600
$ws.Range("K64").Value = 131600
$ws.Range("M64").Value = -131352

$ws.Range("H65").Value = 111113240
$ws.Range("I65").Value = 111113240
$ws.Range("K65").Value = 555566200
$ws.Range("M65").Value = -555563080

$ws.Range("H67").Value = 92248.625
$ws.Range("I67").Value = 131600
$ws.Range("K67").Value = 131600
$ws.Range("M67").Value = -130742

$ws.Range("H76").Value = 4571
$ws.Range("I76").Value = 4353.6665
$ws.Range("J76").Value = 4897
$ws.Range("K76").Value = 4353.6665
$ws.Range("L76").Value = 4897
$ws.Range("M76").Value = -4038.6665
$ws.Range("N76").Value = -5527

$ws.Range("H79").Value = 4571
$ws.Range("I79").Value = 4353.6665
$ws.Range("J79").Value = 4897
$ws.Range("K79").Value = 4353.6665
$ws.Range("L79").Value = 4897
$ws.Range("M79").Value = -3261.6665
$ws.Range("N79").Value = -7081

$ws.Range("H100").Value = 25150396
$ws.Range("I100").Value = 41834996
$ws.Range("K100").Value = 41834996
$ws.Range("M100").Value = -41834455

$ws.Range("H106").Value = 7626.2144
$ws.Range("I106").Value = 8960.777
$ws.Range("K106").Value = 8960.777
$ws.Range("M106").Value = -8329.777

$ws.Range("H111").Value = 1518.2858
$ws.Range("I111").Value = 1599.5
$ws.Range("J111").Value = 1031
$ws.Range("K111").Value = 4798.5
$ws.Range("L111").Value = 3093
$ws.Range("M111").Value = -1731.5
$ws.Range("N111").Value = -9227

$ws.Range("H112").Value = 2830.6667
$ws.Range("I112").Value = 3511.25
$ws.Range("J112").Value = 2583.182
$ws.Range("K112").Value = 10533.75
$ws.Range("L112").Value = 7749.545999999999
$ws.Range("M112").Value = -9425.75
$ws.Range("N112").Value = -9965.545999999998

$ws.Range("H113").Value = 20749.5
$ws.Range("I113").Value = 14999.333
$ws.Range("K113").Value = 14999.333
$ws.Range("M113").Value = -11745.333

$ws.Range("H116").Value = 2354159.8
$ws.Range("I116").Value = 2824291.5
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 2824291.5
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = -2820849.5
$ws.Range("N116").Value = -10384

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 500247
$ws.Range("I2").Value = 494
$ws.Range("J2").Value = 1000000
$ws.Range("K2").Value = 494
$ws.Range("L2").Value = 1000000
$ws.Range("M2").Value = -381
$ws.Range("N2").Value = -1000226

$ws.Range("H4").Value = 98.4
$ws.Range("I4").Value = 98.4
$ws.Range("K4").Value = 98.4
$ws.Range("M4").Value = 17.59999999999999

$ws.Range("H5").Value = 5597.25
$ws.Range("I5").Value = 5597.25
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5597.25
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -5485.25
$ws.Range("N5").ClearContents()

$ws.Range("H26").Value = 7
$ws.Range("I26").Value = 7
$ws.Range("K26").Value = 7
$ws.Range("M26").Value = 323

$ws.Range("H32").Value = 6872.85
$ws.Range("I32").Value = 6738.8037
$ws.Range("J32").Value = 8749.5
$ws.Range("K32").Value = 6738.8037
$ws.Range("L32").Value = 8749.5
$ws.Range("M32").Value = -6451.8037
$ws.Range("N32").Value = -9323.5

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H45").Value = 53798.35
$ws.Range("I45").Value = 94444.73
$ws.Range("K45").Value = 94444.73
$ws.Range("M45").Value = -94067.73

$ws.Range("H61").Value = 9000.464
$ws.Range("I61").Value = 10738.526
$ws.Range("J61").Value = 5331.222
$ws.Range("K61").Value = 10738.526
$ws.Range("L61").Value = 5331.222
$ws.Range("M61").Value = -10526.526
$ws.Range("N61").Value = -5755.222

$ws.Range("H63").Value = 2348.6428
$ws.Range("I63").Value = 2240.5
$ws.Range("K63").Value = 2240.5
$ws.Range("M63").Value = -1554.5

$ws.Range("H66").Value = 2348.6428
$ws.Range("I66").Value = 2240.5
$ws.Range("K66").Value = 11202.5
$ws.Range("M66").Value = -7770.5

$ws.Range("H74").Value = 7338.7144
$ws.Range("I74").Value = 9516.385
$ws.Range("J74").Value = 3800
$ws.Range("K74").Value = 9516.385
$ws.Range("L74").Value = 3800
$ws.Range("M74").Value = -8642.385
$ws.Range("N74").Value = -5548

$ws.Range("H77").Value = 7338.7144
$ws.Range("I77").Value = 9516.385
$ws.Range("J77").Value = 3800
$ws.Range("K77").Value = 47581.925
$ws.Range("L77").Value = 19000
$ws.Range("M77").Value = -43213.925
$ws.Range("N77").Value = -27736

$ws.Range("H88").Value = 55556584
$ws.Range("I88").Value = 537.1429
$ws.Range("J88").Value = 90910430
$ws.Range("K88").Value = 537.1429
$ws.Range("L88").Value = 90910430
$ws.Range("M88").Value = -131.1429000000001
$ws.Range("N88").Value = -90911242

$ws.Range("H91").Value = 55556584
$ws.Range("I91").Value = 537.1429
$ws.Range("J91").Value = 90910430
$ws.Range("K91").Value = 537.1429
$ws.Range("L91").Value = 90910430
$ws.Range("M91").Value = 866.8571
$ws.Range("N91").Value = -90913238

$ws.Range("H97").Value = 8700583
$ws.Range("I97").Value = 7640.5
$ws.Range("J97").Value = 22222938
$ws.Range("K97").Value = 7640.5
$ws.Range("L97").Value = 22222938
$ws.Range("M97").Value = -7144.5
$ws.Range("N97").Value = -22223930

$ws.Range("H110").Value = 2707.8667
$ws.Range("I110").Value = 2010.7273
$ws.Range("K110").Value = 2010.7273
$ws.Range("M110").Value = 34.27269999999999

$ws.Range("H116").Value = 500247
$ws.Range("I116").Value = 494
$ws.Range("J116").Value = 1000000
$ws.Range("K116").Value = 494
$ws.Range("L116").Value = 1000000
$ws.Range("M116").Value = 1800
$ws.Range("N116").Value = -1004588

$ws.Range("H122").Value = 886554.25
$ws.Range("I122").Value = 3633.5833
$ws.Range("K122").Value = 10900.7499
$ws.Range("M122").Value = -8450.749899999999

$ws.Range("H132").Value = 3203.4102
$ws.Range("I132").Value = 2982.2068
$ws.Range("J132").Value = 3844.9
$ws.Range("K132").Value = 8946.6204
$ws.Range("L132").Value = 11534.7
$ws.Range("M132").Value = -6416.6204
$ws.Range("N132").Value = -16594.7

$ws.Range("H136").Value = 9000.464
$ws.Range("I136").Value = 10738.526
$ws.Range("J136").Value = 5331.222
$ws.Range("K136").Value = 32215.578
$ws.Range("L136").Value = 15993.666
$ws.Range("M136").Value = -29665.578
$ws.Range("N136").Value = -21093.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 500247
$ws.Range("I3").Value = 494
$ws.Range("J3").Value = 1000000
$ws.Range("K3").Value = 494
$ws.Range("L3").Value = 1000000
$ws.Range("M3").Value = -380
$ws.Range("N3").Value = -1000228

$ws.Range("H4").Value = 5597.25
$ws.Range("I4").Value = 5597.25
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 5597.25
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -5482.25
$ws.Range("N4").ClearContents()

$ws.Range("H10").Value = 200
$ws.Range("I10").Value = 200
$ws.Range("K10").Value = 200
$ws.Range("M10").Value = -60

$ws.Range("H18").Value = 4500
$ws.Range("J18").Value = 4500
$ws.Range("L18").Value = 4500
$ws.Range("N18").Value = -5558

$ws.Range("H42").Value = 199555
$ws.Range("J42").Value = 199555
$ws.Range("L42").Value = 199555
$ws.Range("N42").Value = -200211

$ws.Range("H94").Value = 8603.058
$ws.Range("I94").Value = 12287.137
$ws.Range("J94").Value = 2368.4614
$ws.Range("K94").Value = 12287.137
$ws.Range("L94").Value = 2368.4614
$ws.Range("M94").Value = -11836.137
$ws.Range("N94").Value = -3270.4614

$ws.Range("H99").Value = 37341.375
$ws.Range("I99").Value = 84440
$ws.Range("K99").Value = 84440
$ws.Range("M99").Value = -82942

$ws.Range("H105").Value = 118210.89
$ws.Range("I105").Value = 204379.6
$ws.Range("K105").Value = 204379.6
$ws.Range("M105").Value = -202632.6

$ws.Range("H107").Value = 2697.111
$ws.Range("I107").Value = 2909.25
$ws.Range("K107").Value = 2909.25
$ws.Range("M107").Value = -989.25

$ws.Range("H134").Value = 5414.9062
$ws.Range("I134").Value = 6304.32
$ws.Range("J134").Value = 2238.4285
$ws.Range("K134").Value = 18912.96
$ws.Range("L134").Value = 6715.2855
$ws.Range("M134").Value = -16377.96
$ws.Range("N134").Value = -11785.2855

$ws.Range("H140").Value = 84213
$ws.Range("J140").Value = 84213
$ws.Range("L140").Value = 84213
$ws.Range("N140").Value = -94573

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 17057.334
$ws.Range("I7").Value = 25242.25
$ws.Range("K7").Value = 25242.25
$ws.Range("M7").Value = -25129.25

$ws.Range("H16").Value = 1502
$ws.Range("I16").Value = 1354.7222
$ws.Range("K16").Value = 1354.7222
$ws.Range("M16").Value = -1067.7222

$ws.Range("H23").Value = 4750
$ws.Range("I23").Value = 4500
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 4500
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = -4260
$ws.Range("N23").Value = -5480

$ws.Range("H27").Value = 4750
$ws.Range("I27").Value = 4500
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 4500
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -4308
$ws.Range("N27").Value = -5384

$ws.Range("H31").Value = 7603.6875
$ws.Range("I31").Value = 8822.211
$ws.Range("J31").Value = 5822.769
$ws.Range("K31").Value = 8822.211
$ws.Range("L31").Value = 5822.769
$ws.Range("M31").Value = -8527.211
$ws.Range("N31").Value = -6412.769

$ws.Range("H34").Value = 7603.6875
$ws.Range("I34").Value = 8822.211
$ws.Range("J34").Value = 5822.769
$ws.Range("K34").Value = 8822.211
$ws.Range("L34").Value = 5822.769
$ws.Range("M34").Value = -8620.211
$ws.Range("N34").Value = -6226.769

$ws.Range("H57").Value = 20222.2
$ws.Range("J57").Value = 20222.2
$ws.Range("L57").Value = 20222.2
$ws.Range("N57").Value = -21342.2

$ws.Range("H58").Value = 2525.5925
$ws.Range("I58").Value = 2426.5
$ws.Range("J58").Value = 2961.6
$ws.Range("K58").Value = 2426.5
$ws.Range("L58").Value = 2961.6
$ws.Range("M58").Value = -2223.5
$ws.Range("N58").Value = -3367.6

$ws.Range("H62").Value = 6313.273
$ws.Range("J62").Value = 6192
$ws.Range("L62").Value = 6192
$ws.Range("N62").Value = -7440

$ws.Range("H65").Value = 6313.273
$ws.Range("J65").Value = 6192
$ws.Range("L65").Value = 30960
$ws.Range("N65").Value = -37200

$ws.Range("H86").Value = 9543.625
$ws.Range("I86").Value = 7391.8335
$ws.Range("K86").Value = 7391.8335
$ws.Range("M86").Value = -6268.8335

$ws.Range("H89").Value = 9543.625
$ws.Range("I89").Value = 7391.8335
$ws.Range("K89").Value = 36959.1675
$ws.Range("M89").Value = -31343.1675

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H99").Value = 197146.27
$ws.Range("I99").Value = 360986.22
$ws.Range("J99").Value = 5999.6665
$ws.Range("K99").Value = 360986.22
$ws.Range("L99").Value = 5999.6665
$ws.Range("M99").Value = -359488.22
$ws.Range("N99").Value = -8995.6665

$ws.Range("H107").Value = 3850.0833
$ws.Range("I107").Value = 4767.607
$ws.Range("K107").Value = 4767.607
$ws.Range("M107").Value = -2847.607

$ws.Range("H109").Value = 37630
$ws.Range("J109").Value = 55000
$ws.Range("L109").Value = 55000
$ws.Range("N109").Value = -57080

$ws.Range("H113").Value = 1502
$ws.Range("I113").Value = 1354.7222
$ws.Range("K113").Value = 1354.7222
$ws.Range("M113").Value = 815.2778000000001

$ws.Range("H122").Value = 10929.571
$ws.Range("I122").Value = 10539.615
$ws.Range("J122").Value = 11563.25
$ws.Range("K122").Value = 31618.845
$ws.Range("L122").Value = 34689.75
$ws.Range("M122").Value = -29168.845
$ws.Range("N122").Value = -39589.75

$ws.Range("H126").Value = 197146.27
$ws.Range("I126").Value = 360986.22
$ws.Range("J126").Value = 5999.6665
$ws.Range("K126").Value = 1082958.66
$ws.Range("L126").Value = 17998.9995
$ws.Range("M126").Value = -1080488.66
$ws.Range("N126").Value = -22938.9995

$ws.Range("H134").Value = 5614.375
$ws.Range("I134").Value = 7657.55
$ws.Range("J134").Value = 2209.0833
$ws.Range("K134").Value = 22972.65
$ws.Range("L134").Value = 6627.249899999999
$ws.Range("M134").Value = -20437.65
$ws.Range("N134").Value = -11697.2499

$ws.Range("H136").Value = 2525.5925
$ws.Range("I136").Value = 2426.5
$ws.Range("J136").Value = 2961.6
$ws.Range("K136").Value = 7279.5
$ws.Range("L136").Value = 8884.8
$ws.Range("M136").Value = -4729.5
$ws.Range("N136").Value = -13984.8

$ws.Range("H138").Value = 75000
$ws.Range("J138").Value = 75000
$ws.Range("L138").Value = 75000
$ws.Range("N138").Value = -85280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 27.1875
$ws.Range("I2").Value = 18.818182
$ws.Range("K2").Value = 112.909092
$ws.Range("M2").Value = 0.09090799999999888

$ws.Range("H5").Value = 625702.25
$ws.Range("J5").Value = 1429917
$ws.Range("L5").Value = 4289751
$ws.Range("N5").Value = -4289975

$ws.Range("H12").Value = 90.210526
$ws.Range("I12").Value = 207.85715
$ws.Range("J12").Value = 21.583334
$ws.Range("K12").Value = 623.5714499999999
$ws.Range("L12").Value = 64.750002
$ws.Range("M12").Value = -450.5714499999999
$ws.Range("N12").Value = -410.750002

$ws.Range("H38").Value = 1365.069
$ws.Range("I38").Value = 208.5
$ws.Range("J38").Value = 1973.7894
$ws.Range("K38").Value = 625.5
$ws.Range("L38").Value = 5921.3682
$ws.Range("M38").Value = -278.5
$ws.Range("N38").Value = -6615.3682

$ws.Range("H40").Value = 69.375
$ws.Range("I40").Value = 28
$ws.Range("J40").Value = 110.75
$ws.Range("K40").Value = 112
$ws.Range("L40").Value = 443
$ws.Range("M40").Value = -43
$ws.Range("N40").Value = -581

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()

$ws.Range("H86").Value = 422.2857
$ws.Range("I86").Value = 407.25
$ws.Range("J86").Value = 442.33334
$ws.Range("K86").Value = 1221.75
$ws.Range("L86").Value = 1327.00002
$ws.Range("M86").Value = -35.75
$ws.Range("N86").Value = -3699.00002

$ws.Range("H87").Value = 17884
$ws.Range("I87").Value = 15909.333
$ws.Range("J87").Value = 18542.223
$ws.Range("K87").Value = 47727.999
$ws.Range("L87").Value = 55626.66900000001
$ws.Range("M87").Value = -46479.999
$ws.Range("N87").Value = -58122.66900000001

$ws.Range("H89").Value = 422.2857
$ws.Range("I89").Value = 407.25
$ws.Range("J89").Value = 442.33334
$ws.Range("K89").Value = 3665.25
$ws.Range("L89").Value = 3981.00006
$ws.Range("M89").Value = 2262.75
$ws.Range("N89").Value = -15837.00006

$ws.Range("H90").Value = 17884
$ws.Range("I90").Value = 15909.333
$ws.Range("J90").Value = 18542.223
$ws.Range("K90").Value = 143183.997
$ws.Range("L90").Value = 166880.007
$ws.Range("M90").Value = -136943.997
$ws.Range("N90").Value = -179360.007

$ws.Range("H107").Value = 1671.1666
$ws.Range("J107").Value = 1782.1818
$ws.Range("L107").Value = 5346.5454
$ws.Range("N107").Value = -9186.545399999999

$ws.Range("H113").Value = 16992.715
$ws.Range("J113").Value = 19658.166
$ws.Range("L113").Value = 58974.49800000001
$ws.Range("N113").Value = -63314.49800000001

$ws.Range("H126").Value = 16527.5
$ws.Range("I126").Value = 3500
$ws.Range("K126").Value = 10500
$ws.Range("M126").Value = -5560

$ws.Range("H131").Value = 1548.73
$ws.Range("I131").Value = 1466.3334
$ws.Range("J131").Value = 1551.2783
$ws.Range("K131").Value = 4399.0002
$ws.Range("L131").Value = 4653.8349
$ws.Range("M131").Value = 640.9997999999996
$ws.Range("N131").Value = -14733.8349

$ws.Range("H135").Value = 625702.25
$ws.Range("J135").Value = 1429917
$ws.Range("L135").Value = 12869253
$ws.Range("N135").Value = -12874323

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 145.91667
$ws.Range("I2").Value = 103.052635
$ws.Range("K2").Value = 103.052635
$ws.Range("M2").Value = 9.947365000000005

$ws.Range("H3").Value = 10910674
$ws.Range("I3").Value = 30000576
$ws.Range("J3").Value = 2158.1428
$ws.Range("K3").Value = 30000576
$ws.Range("L3").Value = 2158.1428
$ws.Range("M3").Value = -30000460
$ws.Range("N3").Value = -2390.1428

$ws.Range("H5").Value = 36000

$ws.Range("H10").Value = 300
$ws.Range("I10").Value = 300
$ws.Range("K10").Value = 300
$ws.Range("M10").Value = -131

$ws.Range("H43").Value = 6000
$ws.Range("I43").Value = 6000
$ws.Range("K43").Value = 6000
$ws.Range("M43").Value = -5849

$ws.Range("H46").Value = 34499
$ws.Range("J46").Value = 34499
$ws.Range("L46").Value = 34499
$ws.Range("N46").Value = -34811

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").ClearContents()

$ws.Range("H57").Value = 49999
$ws.Range("J57").Value = 49999
$ws.Range("L57").Value = 49999
$ws.Range("N57").Value = -51639

$ws.Range("H70").Value = 8721.2
$ws.Range("J70").Value = 10294.167
$ws.Range("L70").Value = 10294.167
$ws.Range("N70").Value = -10834.167

$ws.Range("H73").Value = 8721.2
$ws.Range("J73").Value = 10294.167
$ws.Range("L73").Value = 10294.167
$ws.Range("N73").Value = -12166.167

$ws.Range("H102").Value = 8548.272
$ws.Range("I102").Value = 10453.6
$ws.Range("J102").Value = 4465.4287
$ws.Range("K102").Value = 10453.6
$ws.Range("L102").Value = 4465.4287
$ws.Range("M102").Value = -8831.6
$ws.Range("N102").Value = -7709.4287

$ws.Range("H107").Value = 337.85715
$ws.Range("I107").Value = 371.6
$ws.Range("K107").Value = 371.6
$ws.Range("M107").Value = 1548.4

$ws.Range("H113").Value = 10522.615
$ws.Range("I113").Value = 13653.444
$ws.Range("K113").Value = 13653.444
$ws.Range("M113").Value = -11483.444

$ws.Range("H122").Value = 10626.381
$ws.Range("I122").Value = 8787.182
$ws.Range("J122").Value = 12649.5
$ws.Range("K122").Value = 26361.546
$ws.Range("L122").Value = 37948.5
$ws.Range("M122").Value = -23911.546
$ws.Range("N122").Value = -42848.5

$ws.Range("H126").Value = 10194
$ws.Range("I126").Value = 12952.454
$ws.Range("J126").Value = 8171.1333
$ws.Range("K126").Value = 38857.362
$ws.Range("L126").Value = 24513.3999
$ws.Range("M126").Value = -36387.362
$ws.Range("N126").Value = -29453.3999

$ws.Range("H132").Value = 3642.0393
$ws.Range("I132").Value = 3983.2563
$ws.Range("J132").Value = 2533.0833
$ws.Range("K132").Value = 11949.7689
$ws.Range("L132").Value = 7599.249899999999
$ws.Range("M132").Value = -9419.7689
$ws.Range("N132").Value = -12659.2499

$ws.Range("H134").Value = 93548.664
$ws.Range("J134").Value = 93548.664
$ws.Range("L134").Value = 280645.992
$ws.Range("N134").Value = -285715.992

$ws.Range("H136").Value = 51665
$ws.Range("J136").Value = 51665
$ws.Range("L136").Value = 154995
$ws.Range("N136").Value = -160095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 41243.25
$ws.Range("I7").Value = 74486.836
$ws.Range("J7").Value = 7999.6665
$ws.Range("K7").Value = 74486.836
$ws.Range("L7").Value = 7999.6665
$ws.Range("M7").Value = -74374.836
$ws.Range("N7").Value = -8223.6665

$ws.Range("H16").Value = 3609
$ws.Range("I16").Value = 3560.5
$ws.Range("J16").Value = 3900
$ws.Range("K16").Value = 3560.5
$ws.Range("L16").Value = 3900
$ws.Range("M16").Value = -3390.5
$ws.Range("N16").Value = -4240

$ws.Range("H22").Value = 9490.565
$ws.Range("I22").Value = 14735.071
$ws.Range("J22").Value = 1332.4445
$ws.Range("K22").Value = 14735.071
$ws.Range("L22").Value = 1332.4445
$ws.Range("M22").Value = -14440.071
$ws.Range("N22").Value = -1922.4445

$ws.Range("H27").Value = 9490.565
$ws.Range("I27").Value = 14735.071
$ws.Range("J27").Value = 1332.4445
$ws.Range("K27").Value = 14735.071
$ws.Range("L27").Value = 1332.4445
$ws.Range("M27").Value = -14628.071
$ws.Range("N27").Value = -1546.4445

$ws.Range("H40").Value = 22163.08
$ws.Range("I40").Value = 28216.176
$ws.Range("J40").Value = 9300.25
$ws.Range("K40").Value = 28216.176
$ws.Range("L40").Value = 9300.25
$ws.Range("M40").Value = -28080.176
$ws.Range("N40").Value = -9572.25

$ws.Range("H46").Value = 2324918
$ws.Range("J46").Value = 4102188.8
$ws.Range("L46").Value = 4102188.8
$ws.Range("N46").Value = -4102564.8

$ws.Range("H61").Value = 3824.25
$ws.Range("I61").Value = 825.6667
$ws.Range("J61").Value = 12820
$ws.Range("K61").Value = 825.6667
$ws.Range("L61").Value = 12820
$ws.Range("M61").Value = -623.6667
$ws.Range("N61").Value = -13224

$ws.Range("H68").Value = 6833
$ws.Range("I68").Value = 2549.5
$ws.Range("J68").Value = 15400
$ws.Range("K68").Value = 2549.5
$ws.Range("L68").Value = 15400
$ws.Range("M68").Value = -1800.5
$ws.Range("N68").Value = -16898

$ws.Range("H71").Value = 6833
$ws.Range("I71").Value = 2549.5
$ws.Range("J71").Value = 15400
$ws.Range("K71").Value = 12747.5
$ws.Range("L71").Value = 77000
$ws.Range("M71").Value = -9003.5
$ws.Range("N71").Value = -84488

$ws.Range("H100").Value = 11207.417
$ws.Range("I100").Value = 13498.429
$ws.Range("K100").Value = 13498.429
$ws.Range("M100").Value = -12957.429

$ws.Range("H113").Value = 3824.25
$ws.Range("I113").Value = 825.6667
$ws.Range("J113").Value = 12820
$ws.Range("K113").Value = 825.6667
$ws.Range("L113").Value = 12820
$ws.Range("M113").Value = 1344.3333
$ws.Range("N113").Value = -17160

$ws.Range("H126").Value = 41243.25
$ws.Range("I126").Value = 74486.836
$ws.Range("J126").Value = 7999.6665
$ws.Range("K126").Value = 223460.508
$ws.Range("L126").Value = 23998.9995
$ws.Range("M126").Value = -220990.508
$ws.Range("N126").Value = -28938.9995

$ws.Range("H132").Value = 1068128.8
$ws.Range("I132").Value = 2129915.2
$ws.Range("J132").Value = 6342.143
$ws.Range("K132").Value = 6389745.600000001
$ws.Range("L132").Value = 19026.429
$ws.Range("M132").Value = -6387215.600000001
$ws.Range("N132").Value = -24086.429

$ws.Range("H136").Value = 4285.926
$ws.Range("I136").Value = 2835.3845
$ws.Range("J136").Value = 5632.857
$ws.Range("K136").Value = 8506.1535
$ws.Range("L136").Value = 16898.571
$ws.Range("M136").Value = -5956.1535
$ws.Range("N136").Value = -21998.571

$ws.Range("H140").Value = 84197.8
$ws.Range("J140").Value = 103747.25
$ws.Range("L140").Value = 103747.25
$ws.Range("N140").Value = -114107.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 250
$ws.Range("I6").Value = 250
$ws.Range("K6").Value = 250
$ws.Range("M6").Value = -135

$ws.Range("H18").Value = 1000
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("N18").Value = -1346

$ws.Range("H33").Value = 14000
$ws.Range("J33").Value = 14000
$ws.Range("L33").Value = 14000
$ws.Range("N33").Value = -14500

$ws.Range("H36").Value = 14000
$ws.Range("J36").Value = 14000
$ws.Range("L36").Value = 14000
$ws.Range("N36").Value = -14500

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H100").Value = 45840.23
$ws.Range("I100").Value = 30527.111
$ws.Range("J100").Value = 80294.75
$ws.Range("K100").Value = 61054.222
$ws.Range("L100").Value = 160589.5
$ws.Range("M100").Value = -60513.222
$ws.Range("N100").Value = -161671.5

$ws.Range("H107").Value = 27006
$ws.Range("I107").Value = 2188.3635
$ws.Range("K107").Value = 6565.0905
$ws.Range("M107").Value = -4645.0905

$ws.Range("H112").Value = 45833
$ws.Range("J112").Value = 45833
$ws.Range("L112").Value = 45833
$ws.Range("N112").Value = -48787

$ws.Range("H113").Value = 1549.1428
$ws.Range("I113").Value = 793.4706
$ws.Range("K113").Value = 2380.4118
$ws.Range("M113").Value = -210.4117999999999

$ws.Range("H115").Value = 20000
$ws.Range("J115").Value = 20000
$ws.Range("L115").Value = 20000
$ws.Range("N115").Value = -23134

$ws.Range("H122").Value = 4121.3403
$ws.Range("J122").Value = 7858.5
$ws.Range("L122").Value = 23575.5
$ws.Range("N122").Value = -28475.5

$ws.Range("H126").Value = 47942.445
$ws.Range("I126").Value = 81896.6
$ws.Range("K126").Value = 245689.8
$ws.Range("M126").Value = -243219.8

$ws.Range("H132").Value = 20186.174
$ws.Range("I132").Value = 37739.91
$ws.Range("J132").Value = 4095.25
$ws.Range("K132").Value = 113219.73
$ws.Range("L132").Value = 12285.75
$ws.Range("M132").Value = -110689.73
$ws.Range("N132").Value = -17345.75

$ws.Range("H136").Value = 683066.3
$ws.Range("I136").Value = 973745.3
$ws.Range("J136").Value = 18657.143
$ws.Range("K136").Value = 2921235.9
$ws.Range("L136").Value = 55971.429
$ws.Range("M136").Value = -2918685.9
$ws.Range("N136").Value = -61071.429
